$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Authors cell: "Aran E." -> append ", Berlanga A." as a new run ---
$cellAuthor = $t.Cell(3, 2)
$rAuthor = $cellAuthor.Range
$insAuthor = $d.Range($rAuthor.End - 1, $rAuthor.End - 1)
$insAuthor.InsertAfter(", Berlanga A.")

# --- 2. Fecha de creacion cell: "20-FEBRERO-2024" -> append "; 06-MARZO-2024" as a new run ---
$cellFecha = $t.Cell(4, 2)
$rFecha = $cellFecha.Range
$insFecha = $d.Range($rFecha.End - 1, $rFecha.End - 1)
$insFecha.InsertAfter("; 06-MARZO-2024")

# --- 3. Flujo alternativo cell: replace text and turn into a numbered list item ---
$cellFlujo = $t.Cell(8, 2)
$rFlujo = $cellFlujo.Range
$textRangeFlujo = $d.Range($rFlujo.Start, $rFlujo.End - 1)
$textRangeFlujo.Text = "Si no consigue; Se busca otro proveedor o se queda en espera."

$pFlujo = $cellFlujo.Range.Paragraphs.Item(1)
$pFlujo.Range.Style = "Prrafodelista"
$pFlujo.Range.ListFormat.ApplyNumberDefault()

Write-Host "done"
